$wb = $excel.ActiveWorkbook

# Remove the old lookup "UserDetails" sheet (Name/Country data); the surviving
# sheet becomes the renamed "Sheet1" (firstName/lastName/... test-data sheet).
[void]$wb.Worksheets.Item("UserDetails").Delete()

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "UserDetails"
$ws1.Activate()

# Simulate the repeated CreateData()/UpdatedDataInExcel() runs that generated
# the final row of test data (firstName, lastName, userName, password, email,
# mobileNumber) written in column order: first, last, password, mobile, user, email.
$ws1.Range("A2").Value = "Saul"
$ws1.Range("B2").Value = "Koch"
$ws1.Range("D2").Value = "6g4488swo2p5"
$ws1.Range("F2").Value = "250-664-9473"
$ws1.Range("C2").Value = "SaulKoch56660"
$ws1.Range("E2").Value = "drew.kiehn@hotmail.com"

[void]$ws1.Range("F2").Select()
